# Correction in sa algorithm and 746 logs
# Update the "Fitness" (column C) values for generations 0-57 (rows 2-59)
# on the active worksheet of run_14.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFitness = @(
    8184, 8184, 8184, 8184, 8184, 8184, 8184, 8184, 8184, 8184, 8184, 8184, 8184,
    8159,
    7765, 7765, 7765, 7765, 7765, 7765, 7765,
    7310, 7310, 7310, 7310, 7310, 7310, 7310, 7310, 7310, 7310, 7310, 7310,
    7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293,
    7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293, 7293
)

$startRow = 2
for ($i = 0; $i -lt $newFitness.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newFitness[$i]
}
